$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the ranges we touch so numeric-looking
# strings (e.g. "0.9998") are preserved as literal text, not converted
# to numbers by Excel's auto-detection.
$ws.Range("B9:E50").NumberFormat = "@"
$ws.Range("D2:E8").NumberFormat = "@"
$ws.Range("D51:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.899.48'
$ws.Range("E2").Value = '  +4.48%  '
$ws.Range("D3").Value = '1.877.74'
$ws.Range("E3").Value = '  +3.57%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '278.97'
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").Value = '0.5278'
$ws.Range("E7").Value = '  +4.24%  '
$ws.Range("D8").Value = '0.3448'
$ws.Range("E8").Value = '  -1.77%  '
$ws.Range("D9").Value = '0.06950'
$ws.Range("E9").Value = '  +4.12%  '
$ws.Range("D10").Value = '20.15'
$ws.Range("E10").Value = '  +0.73%  '
$ws.Range("D11").Value = '0.8079'
$ws.Range("E11").Value = '  -2.85%  '
$ws.Range("D12").Value = '0.07869'
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").Value = '1.853.58'
$ws.Range("E13").Value = '  +2.32%  '
$ws.Range("D14").Value = '90.11'
$ws.Range("E14").Value = '  +3.02%  '
$ws.Range("D15").Value = '5.168'
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("D16").Value = '14.59'
$ws.Range("E16").Value = '  +4.20%  '
$ws.Range("D17").Value = '0.9992'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").Value = '0.000008066'
$ws.Range("E18").Value = '  +0.31%  '
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = '26.944.08'
$ws.Range("E20").Value = '  +4.46%  '
$ws.Range("D21").Value = '2.105.34'
$ws.Range("E21").Value = '  +3.43%  '
$ws.Range("D22").Value = '4.753'
$ws.Range("E22").Value = '  +0.77%  '
$ws.Range("D23").Value = '10.03'
$ws.Range("E23").Value = '  +0.33%  '
$ws.Range("D24").Value = '6.194'
$ws.Range("E24").Value = '  +2.37%  '
$ws.Range("D25").Value = '2.348'
$ws.Range("E25").Value = '  +7.82%  '
$ws.Range("D26").Value = '146.59'
$ws.Range("E26").Value = '  +3.12%  '
$ws.Range("D27").Value = '17.38'
$ws.Range("E27").Value = '  +2.04%  '
$ws.Range("D28").Value = '1.658'
$ws.Range("E28").Value = '  -0.47%  '
$ws.Range("D29").Value = '113.94'
$ws.Range("E29").Value = '  +4.08%  '
$ws.Range("D30").Value = '4.382'
$ws.Range("E30").Value = '  +1.32%  '
$ws.Range("D31").Value = '4.326'
$ws.Range("E31").Value = '  +2.32%  '
$ws.Range("D32").Value = '0.08924'
$ws.Range("E32").Value = '  +1.41%  '
$ws.Range("D33").Value = '0.04948'
$ws.Range("E33").Value = '  +1.64%  '
$ws.Range("D34").Value = '1.175'
$ws.Range("E34").Value = '  +3.39%  '
$ws.Range("D35").Value = '0.7365'
$ws.Range("E35").Value = '  +1.34%  '
$ws.Range("D36").Value = '2.893'
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("D37").Value = '3.283'
$ws.Range("E37").Value = '  +4.47%  '
$ws.Range("D38").Value = '2.408'
$ws.Range("E38").Value = '  +6.01%  '
$ws.Range("D39").Value = '0.01852'
$ws.Range("E39").Value = '  +0.59%  '
$ws.Range("D40").Value = '0.5147'
$ws.Range("E40").Value = '  -0.89%  '
$ws.Range("D41").Value = '0.9579'
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("D42").Value = '116.43'
$ws.Range("E42").Value = '  +2.68%  '
$ws.Range("D43").Value = '6.211'
$ws.Range("E43").Value = '  +1.14%  '
$ws.Range("D44").Value = '8.095'
$ws.Range("E44").Value = '  +0.38%  '
$ws.Range("D45").Value = '0.9992'
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("D46").Value = '0.4514'
$ws.Range("E46").Value = '  -1.01%  '
$ws.Range("D47").Value = '0.1344'
$ws.Range("E47").Value = '  -1.01%  '
$ws.Range("D48").Value = '9.410'
$ws.Range("E48").Value = '  +1.23%  '
$ws.Range("D49").Value = '36.60'
$ws.Range("E49").Value = '  +0.75%  '
$ws.Range("D50").Value = '1.507'
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("D51").Value = '0.05950'
$ws.Range("E51").Value = '  +2.03%  '

# Coin name / link columns shift for rows 9-50 (new coin list order)
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("B20").Value = 'WrappedBTC'
$ws.Range("C20").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("B25").Value = 'LidoDAOToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'

# Reset style (drop the temporary text number format) so cells keep
# their original default styling, matching the source workbook.
$ws.Range("B9:E50").Style = "Normal"
$ws.Range("D2:E8").Style = "Normal"
$ws.Range("D51:E51").Style = "Normal"
